$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 8-10 (MuSCs -> MuSCs triple of rows) entirely; data now ends at row 7.
$ws.Rows("8:10").Delete()

# Updated data (columns A:T) for rows 2-7 reflecting the new TPM-derived values.
$data = @(
  @("ECs","Icam2","Itgb2","ECs",3,1,16.782487,50.347461,0.9252099721531751,0.9252099721531751,2,0.6666666666666666,0.05619066666666667,0.168572,0.3931387525216601,0.39313875252166,0.9430191328546667,8.487172195692001,0.3637358942728992,0.3637358942728991),
  @("ECs","Icam2","Itgb2","FAPs",3,1,16.782487,50.347461,0.9252099721531751,0.9252099721531751,2,0.6666666666666666,0.08673766666666667,0.260213,0.60686124747834,0.60686124747834,1.455673763243667,13.101063869193,0.5614740778802761,0.5614740778802761),
  @("FAPs","Icam2","Itgb2","ECs",3,1,1.091026,3.273078,0.0601477084462148,0.0601477084462148,2,0.6666666666666666,0.05619066666666667,0.168572,0.3931387525216601,0.39313875252166,0.06130547829066667,0.551749304616,0.02364639506558141,0.0236463950655814),
  @("FAPs","Icam2","Itgb2","FAPs",3,1,1.091026,3.273078,0.0601477084462148,0.0601477084462148,2,0.6666666666666666,0.08673766666666667,0.260213,0.60686124747834,0.60686124747834,0.09463304951266667,0.8516974456140001,0.03650131338063341,0.0365013133806334),
  @("MuSCs","Icam2","Itgb2","ECs",2,0.6666666666666666,0.2655986666666667,0.7967960000000001,0.01464231940061012,0.01464231940061012,2,0.6666666666666666,0.05619066666666667,0.168572,0.3931387525216601,0.39313875252166,0.01492416614577778,0.134317495312,0.005756463183179565,0.005756463183179563),
  @("MuSCs","Icam2","Itgb2","FAPs",2,0.6666666666666666,0.2655986666666667,0.7967960000000001,0.01464231940061012,0.01464231940061012,2,0.6666666666666666,0.08673766666666667,0.260213,0.60686124747834,0.60686124747834,0.02303740861644445,0.207336677548,0.008885856217430559,0.008885856217430557)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $r = $i + 2
  $row = $data[$i]
  for ($c = 0; $c -lt $row.Length; $c++) {
    $ws.Cells.Item($r, $c + 1).Value = $row[$c]
  }
}

$wb.Save()
